# Commit: "Added filtering options for the Component Analysis"
#
# The underlying edit clears the trailing-quarter error columns from the
# naive_err / AVERAGE_1_9_qoq_errors staircase table: for each data row
# (rows 2-44), everything to the right of the "current" filter window is
# cleared back to blank, leaving the staircase taper that continues
# unmodified into rows 45-53. Cell contents are cleared (not the rows/
# columns themselves), so row spans and sheet dimension stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ranges = @(
    "G2:K2",
    "I3:K3",
    "G4:K4",
    "I5:K5",
    "G6:K6",
    "I7:K7",
    "G8:K8",
    "I9:K9",
    "G10:K10",
    "I11:K11",
    "G12:K12",
    "I13:K13",
    "G14:K14",
    "I15:K15",
    "G16:K16",
    "I17:K17",
    "K18",
    "I19:K19",
    "K20",
    "I21:K21",
    "K22",
    "J23:K23",
    "I24:K24",
    "K26",
    "J27:K27",
    "I28:K28",
    "K30",
    "J31:K31",
    "I32:K32",
    "K34",
    "J35:K35",
    "I36:K36",
    "K38",
    "J39:K39",
    "I40:K40",
    "K42",
    "J43:K43",
    "I44:J44"
)

foreach ($addr in $ranges) {
    $ws.Range($addr).ClearContents()
}
